$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.127.27"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.792.23"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.64"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0716"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "2.049.35"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "1.796.48"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.625"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "34.076.27"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("D20").Value = "0.0₃0792"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("D35").Value = "1.399.69"
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.04%  "
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.921"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.20%  "
$ws.Range("E45").Value = "  +3.48%  "
$ws.Range("D46").Value = "0.0₆0137"
$ws.Range("E46").Value = "  +10.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "109.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").Value = "1.949.65"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.23%  "
